$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.805.73"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "3.096.70"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'541.34"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'137.79"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.091.41"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "'6.44"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'0.0000227"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("D14").Value = "'34.73"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "3.589.09"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "63.795.28"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "3.089.50"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'6.69"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "'482.64"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'13.43"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'7.12"
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("D24").Value = "'79.35"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("D25").Value = "'12.24"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'8.10"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'26.32"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'1.90"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "'57.58"
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "'498.87"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "'5.35"
$ws.Range("E36").Value = "  +6.06%  "
$ws.Range("D37").Value = "'6.01"
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").Value = "3.257.57"
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("D39").Value = "'0.0404"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").Value = "'0.0796"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "'8.13"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").Value = "'0.255"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D46").Value = "'123.39"
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("D47").Value = "'2.04"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "0.0₃0530"
$ws.Range("E48").Value = "  +9.18%  "
$ws.Range("D49").Value = "'24.64"
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "'2.42"
$ws.Range("E51").Value = "  +4.22%  "
